# Apply the "Penalty Reward System" forecast-shift edit:
#  - Sheet "Forecast Comparison": shift Week_Start_Date (col B) forward by one
#    week for rows 2-17, and update MyForecast (col D) to the new values.
#  - Sheet "Summary": refresh the derived statistics (col B) to match the
#    updated forecast column.

$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# New Week_Start_Date values (col B), rows 2-17
$newDates = @(
    "2025-01-12",
    "2025-01-19",
    "2025-01-26",
    "2025-02-02",
    "2025-02-09",
    "2025-02-16",
    "2025-02-23",
    "2025-03-02",
    "2025-03-09",
    "2025-03-16",
    "2025-03-23",
    "2025-03-30",
    "2025-04-06",
    "2025-04-13",
    "2025-04-20",
    "2025-04-27"
)

# New MyForecast values (col D), rows 2-17
$newForecast = @(107, 102, 96, 96, 104, 110, 108, 102, 100, 107, 99, 99, 107, 103, 96, 97)

# Force column B to text so the date-like strings aren't silently converted to
# Excel date serial numbers; restore default formatting afterwards so we
# don't leave stray number-format styling behind.
$dateRange = $wsForecast.Range("B2:B17")
$dateRange.NumberFormat = "@"

for ($i = 0; $i -lt $newDates.Length; $i++) {
    $row = $i + 2
    $wsForecast.Cells.Item($row, 2).Value = $newDates[$i]
    $wsForecast.Cells.Item($row, 4).Value = $newForecast[$i]
}

$dateRange.ClearFormats()

# ---- Summary sheet -------------------------------------------------------
# Every value in column B of the Summary sheet is stored as text (even the
# purely numeric-looking ones), so force text formatting before writing to
# avoid Excel silently re-typing them as numbers/dates, then clear the
# formatting afterwards to leave no stray style behind.

$summaryRange = $wsSummary.Range("B2:B15")
$summaryRange.NumberFormat = "@"

$wsSummary.Range("B2").Value  = "2024-05-12 to 2025-01-05"
$wsSummary.Range("B4").Value  = "108"
$wsSummary.Range("B5").Value  = "42"
$wsSummary.Range("B6").Value  = "38"
$wsSummary.Range("B7").Value  = "34"
$wsSummary.Range("B8").Value  = "1458 units"
$wsSummary.Range("B9").Value  = "1633"
$wsSummary.Range("B10").Value = "825"
$wsSummary.Range("B11").Value = "401"
$wsSummary.Range("B12").Value = "110"
$wsSummary.Range("B13").Value = "2025-02-16"
$wsSummary.Range("B14").Value = "96"
$wsSummary.Range("B15").Value = "2025-01-26"

$summaryRange.ClearFormats()
